$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.276276707649231
$ws.Range("B1").Value = 2.390237331390381
$ws.Range("D1").Value = 1.376641392707825
$ws.Range("E1").Value = 0.85772705078125
